# Update gh-pages output data (generated at 456a3b4)
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5644
$ws1.Range("F5").Value = 312
$ws1.Range("F6").Value = 838
$ws1.Range("F7").Value = 57
$ws1.Range("F8").Value = 376
$ws1.Range("F11").Value = 21

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 50

# Sheet "全部类型" (all types combined)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5644
$ws4.Range("F5").Value = 312
$ws4.Range("F6").Value = 838
$ws4.Range("F7").Value = 57
$ws4.Range("F8").Value = 50
$ws4.Range("F9").Value = 376
$ws4.Range("F12").Value = 21
